$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 126.5
$ws.Range("I6").Value = 126.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 379.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -267.5
$ws.Range("N6").ClearContents()

$ws.Range("H33").Value = 147.53847
$ws.Range("I33").Value = 112.22222
$ws.Range("K33").Value = 112.22222
$ws.Range("M33").Value = 116.77778

$ws.Range("H38").Value = 395.46155
$ws.Range("J38").Value = 949.25
$ws.Range("L38").Value = 2847.75
$ws.Range("N38").Value = -3591.75

$ws.Range("H39").Value = 477.8
$ws.Range("I39").Value = 390
$ws.Range("J39").Value = 499.75
$ws.Range("K39").Value = 1170
$ws.Range("L39").Value = 1499.25
$ws.Range("M39").Value = -874
$ws.Range("N39").Value = -2091.25

$ws.Range("H129").Value = 928.09
$ws.Range("I129").Value = 432.2
$ws.Range("J129").Value = 954.18945
$ws.Range("K129").Value = 1296.6
$ws.Range("L129").Value = 2862.56835
$ws.Range("M129").Value = 3703.4
$ws.Range("N129").Value = -12862.56835

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9331.125
$ws.Range("I61").Value = 6429.6387
$ws.Range("J61").Value = 18035.584
$ws.Range("K61").Value = 6429.6387
$ws.Range("L61").Value = 18035.584
$ws.Range("M61").Value = -6217.6387
$ws.Range("N61").Value = -18459.584

$ws.Range("H74").Value = 13293.8
$ws.Range("I74").Value = 3844.4285
$ws.Range("J74").Value = 35342.332
$ws.Range("K74").Value = 3844.4285
$ws.Range("L74").Value = 35342.332
$ws.Range("M74").Value = -2970.4285
$ws.Range("N74").Value = -37090.332

$ws.Range("H77").Value = 13293.8
$ws.Range("I77").Value = 3844.4285
$ws.Range("J77").Value = 35342.332
$ws.Range("K77").Value = 19222.1425
$ws.Range("L77").Value = 176711.66
$ws.Range("M77").Value = -14854.1425
$ws.Range("N77").Value = -185447.66

$ws.Range("H136").Value = 9331.125
$ws.Range("I136").Value = 6429.6387
$ws.Range("J136").Value = 18035.584
$ws.Range("K136").Value = 19288.9161
$ws.Range("L136").Value = 54106.75199999999
$ws.Range("M136").Value = -16738.9161
$ws.Range("N136").Value = -59206.75199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1071.0714
$ws.Range("I99").Value = 1071.0714
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1071.0714
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 426.9286
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3222.6333
$ws.Range("I134").Value = 2536.15
$ws.Range("K134").Value = 7608.450000000001
$ws.Range("M134").Value = -5073.450000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5182.864
$ws.Range("I3").Value = 3224.9
$ws.Range("J3").Value = 6814.5
$ws.Range("K3").Value = 9674.700000000001
$ws.Range("L3").Value = 20443.5
$ws.Range("M3").Value = -9562.700000000001
$ws.Range("N3").Value = -20667.5

$ws.Range("H10").Value = 799.3333
$ws.Range("I10").Value = 799.3333
$ws.Range("K10").Value = 2397.9999
$ws.Range("M10").Value = -2258.9999

$ws.Range("H11").Value = 286
$ws.Range("I11").Value = 256.33334
$ws.Range("J11").Value = 375
$ws.Range("K11").Value = 769.0000200000001
$ws.Range("L11").Value = 1125
$ws.Range("M11").Value = -629.0000200000001
$ws.Range("N11").Value = -1405

$ws.Range("H22").Value = 1670.8334
$ws.Range("I22").Value = 750
$ws.Range("J22").Value = 1977.7778
$ws.Range("K22").Value = 2250
$ws.Range("L22").Value = 5933.3334
$ws.Range("M22").Value = -2081
$ws.Range("N22").Value = -6271.3334

$ws.Range("H25").Value = 1575.4286
$ws.Range("I25").Value = 405.6
$ws.Range("J25").Value = 4500
$ws.Range("K25").Value = 1216.8
$ws.Range("L25").Value = 13500
$ws.Range("M25").Value = -1047.8
$ws.Range("N25").Value = -13838

$ws.Range("H27").Value = 1670.8334
$ws.Range("I27").Value = 750
$ws.Range("J27").Value = 1977.7778
$ws.Range("K27").Value = 2250
$ws.Range("L27").Value = 5933.3334
$ws.Range("M27").Value = -2148
$ws.Range("N27").Value = -6137.3334

$ws.Range("H30").Value = 1575.4286
$ws.Range("I30").Value = 405.6
$ws.Range("J30").Value = 4500
$ws.Range("K30").Value = 1216.8
$ws.Range("L30").Value = 13500
$ws.Range("M30").Value = -1114.8
$ws.Range("N30").Value = -13704

$ws.Range("H46").Value = 2822
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2822
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 8466
$ws.Range("N46").Value = -8648
$ws.Range("M46").ClearContents()

$ws.Range("H47").Value = 993
$ws.Range("I47").Value = 903
$ws.Range("J47").Value = 1011
$ws.Range("K47").Value = 2709
$ws.Range("L47").Value = 3033
$ws.Range("M47").Value = -2278
$ws.Range("N47").Value = -3895

$ws.Range("H81").Value = 4441.9287
$ws.Range("I81").Value = 4000
$ws.Range("J81").Value = 4475.923
$ws.Range("K81").Value = 12000
$ws.Range("L81").Value = 13427.769
$ws.Range("M81").Value = -10877
$ws.Range("N81").Value = -15673.769

$ws.Range("H84").Value = 4441.9287
$ws.Range("I84").Value = 4000
$ws.Range("J84").Value = 4475.923
$ws.Range("K84").Value = 36000
$ws.Range("L84").Value = 40283.307
$ws.Range("M84").Value = -30384
$ws.Range("N84").Value = -51515.307

$ws.Range("H104").Value = 2254.5
$ws.Range("J104").Value = 2859.3333
$ws.Range("L104").Value = 8577.999899999999
$ws.Range("N104").Value = -13819.9999

$ws.Range("H125").Value = 4396.9287
$ws.Range("J125").Value = 4791.4165
$ws.Range("L125").Value = 14374.2495
$ws.Range("N125").Value = -24214.2495

$ws.Range("H134").Value = 2991.5806
$ws.Range("I134").Value = 2178.1365
$ws.Range("K134").Value = 6534.4095
$ws.Range("M134").Value = -1464.4095

$ws.Range("H139").Value = 1906065.4
$ws.Range("I139").Value = 3202933.5
$ws.Range("J139").Value = 3992.1333
$ws.Range("K139").Value = 9608800.5
$ws.Range("L139").Value = 11976.3999
$ws.Range("M139").Value = -9603660.5
$ws.Range("N139").Value = -22256.3999

$ws.Range("H140").Value = 2433.2307
$ws.Range("I140").Value = 1910.625
$ws.Range("J140").Value = 3269.4
$ws.Range("K140").Value = 5731.875
$ws.Range("L140").Value = 9808.200000000001
$ws.Range("M140").Value = -551.875
$ws.Range("N140").Value = -20168.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 54718.75
$ws.Range("J135").Value = 54718.75
$ws.Range("L135").Value = 54718.75
$ws.Range("N135").Value = -64858.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2000
$ws.Range("I100").Value = 2000
$ws.Range("K100").Value = 2000
$ws.Range("M100").Value = -1459

$ws.Range("H136").Value = 5108.4053
$ws.Range("I136").Value = 2456.1667
$ws.Range("J136").Value = 7621.0527
$ws.Range("K136").Value = 7368.500100000001
$ws.Range("L136").Value = 22863.1581
$ws.Range("M136").Value = -4818.500100000001
$ws.Range("N136").Value = -27963.1581

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5460.0576
$ws.Range("I136").Value = 2748.138
$ws.Range("J136").Value = 8879.434999999999
$ws.Range("K136").Value = 8244.414000000001
$ws.Range("L136").Value = 26638.305
$ws.Range("M136").Value = -5694.414000000001
$ws.Range("N136").Value = -31738.305
